$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 17 to make room for the new "MuSCs" sending-cluster block
# (shifts the existing "Resolving-Mac" block down from rows 17-21 to rows 22-26)
$ws.Range("A17:A21").EntireRow.Insert()

$data = @(
    @("ECs", "Il23a", "Il12rb1", "ECs", 3, 1, 1.270105666666667, 3.810317, 0.220415243968963, 0.220415243968963, 1, 0.3333333333333333, 0.2848286666666667, 0.854486, 0.08022967564521397, 0.08022967564521397, 0.3617625035624444, 3.255862532062, 0.0176838435308906, 0.0176838435308906),
    @("ECs", "Il23a", "Il12rb1", "FAPs", 3, 1, 1.270105666666667, 3.810317, 0.220415243968963, 0.220415243968963, 3, 1, 0.9694063333333333, 2.908219, 0.2730594381596053, 0.2730594381596053, 1.231248477269222, 11.081236295423, 0.06018646267997736, 0.06018646267997736),
    @("ECs", "Il23a", "Il12rb1", "Inflammatory-Mac", 3, 1, 1.270105666666667, 3.810317, 0.220415243968963, 0.220415243968963, 3, 1, 1.351483666666667, 4.054451, 0.3806818235166093, 0.3806818235166093, 1.716527063440778, 15.448743570967, 0.08390807700496314, 0.08390807700496315),
    @("ECs", "Il23a", "Il12rb1", "MuSCs", 3, 1, 1.270105666666667, 3.810317, 0.220415243968963, 0.220415243968963, 3, 1, 0.23571, 0.70713, 0.066394078474077, 0.066394078474077, 0.29937660669, 2.69438946021, 0.01463426700495815, 0.01463426700495815),
    @("ECs", "Il23a", "Il12rb1", "Resolving-Mac", 3, 1, 1.270105666666667, 3.810317, 0.220415243968963, 0.220415243968963, 3, 1, 0.7087373333333332, 2.126212, 0.1996349842044944, 0.1996349842044945, 0.9001713032448887, 8.101541729204, 0.04400259374817371, 0.04400259374817372),
    @("FAPs", "Il23a", "Il12rb1", "ECs", 3, 1, 1.653878333333333, 4.961635, 0.287015486903044, 0.287015486903044, 1, 0.3333333333333333, 0.2848286666666667, 0.854486, 0.08022967564521397, 0.08022967564521397, 0.4710719605122222, 4.23964764461, 0.02302715941938437, 0.02302715941938437),
    @("FAPs", "Il23a", "Il12rb1", "FAPs", 3, 1, 1.653878333333333, 4.961635, 0.287015486903044, 0.287015486903044, 3, 1, 0.9694063333333333, 2.908219, 0.2730594381596053, 0.2730594381596053, 1.603280130896111, 14.429521178065, 0.07837228759685075, 0.07837228759685075),
    @("FAPs", "Il23a", "Il12rb1", "Inflammatory-Mac", 3, 1, 1.653878333333333, 4.961635, 0.287015486903044, 0.287015486903044, 3, 1, 1.351483666666667, 4.054451, 0.3806818235166093, 0.3806818235166093, 2.235189554153889, 20.116705987385, 0.1092615789317583, 0.1092615789317583),
    @("FAPs", "Il23a", "Il12rb1", "MuSCs", 3, 1, 1.653878333333333, 4.961635, 0.287015486903044, 0.287015486903044, 3, 1, 0.23571, 0.70713, 0.066394078474077, 0.066394078474077, 0.38983566195, 3.50852095755, 0.01905612876071612, 0.01905612876071612),
    @("FAPs", "Il23a", "Il12rb1", "Resolving-Mac", 3, 1, 1.653878333333333, 4.961635, 0.287015486903044, 0.287015486903044, 3, 1, 0.7087373333333332, 2.126212, 0.1996349842044944, 0.1996349842044945, 1.172165319624444, 10.54948787662, 0.05729833219433446, 0.05729833219433447),
    @("Inflammatory-Mac", "Il23a", "Il12rb1", "ECs", 2, 0.6666666666666666, 1.447887666666667, 4.343663, 0.2512676871409801, 0.2512676871409801, 1, 0.3333333333333333, 0.2848286666666667, 0.854486, 0.08022967564521397, 0.08022967564521397, 0.4123999135797778, 3.711599222218, 0.02015912503944393, 0.02015912503944393),
    @("Inflammatory-Mac", "Il23a", "Il12rb1", "FAPs", 2, 0.6666666666666666, 1.447887666666667, 4.343663, 0.2512676871409801, 0.2512676871409801, 3, 1, 0.9694063333333333, 2.908219, 0.2730594381596053, 0.2730594381596053, 1.403591474021889, 12.632323266197, 0.0686110134783795, 0.0686110134783795),
    @("Inflammatory-Mac", "Il23a", "Il12rb1", "Inflammatory-Mac", 2, 0.6666666666666666, 1.447887666666667, 4.343663, 0.2512676871409801, 0.2512676871409801, 3, 1, 1.351483666666667, 4.054451, 0.3806818235166093, 0.3806818235166093, 1.956796532668111, 17.611168794013, 0.09565304133162918, 0.09565304133162919),
    @("Inflammatory-Mac", "Il23a", "Il12rb1", "MuSCs", 2, 0.6666666666666666, 1.447887666666667, 4.343663, 0.2512676871409801, 0.2512676871409801, 3, 1, 0.23571, 0.70713, 0.066394078474077, 0.066394078474077, 0.34128160191, 3.071534417190001, 0.01668268653803806, 0.01668268653803806),
    @("Inflammatory-Mac", "Il23a", "Il12rb1", "Resolving-Mac", 2, 0.6666666666666666, 1.447887666666667, 4.343663, 0.2512676871409801, 0.2512676871409801, 3, 1, 0.7087373333333332, 2.126212, 0.1996349842044944, 0.1996349842044945, 1.026172043839555, 9.235548394556, 0.05016182075348941, 0.05016182075348941),
    @("MuSCs", "Il23a", "Il12rb1", "ECs", 2, 0.6666666666666666, 0.324919, 0.974757, 0.05638672634467276, 0.05638672634467275, 1, 0.3333333333333333, 0.2848286666666667, 0.854486, 0.08022967564521397, 0.08022967564521397, 0.09254624554466667, 0.832916209902, 0.004523888765328537, 0.004523888765328536),
    @("MuSCs", "Il23a", "Il12rb1", "FAPs", 2, 0.6666666666666666, 0.324919, 0.974757, 0.05638672634467276, 0.05638672634467275, 3, 1, 0.9694063333333333, 2.908219, 0.2730594381596053, 0.2730594381596053, 0.3149785364203333, 2.834806827783, 0.01539692781533576, 0.01539692781533576),
    @("MuSCs", "Il23a", "Il12rb1", "Inflammatory-Mac", 2, 0.6666666666666666, 0.324919, 0.974757, 0.05638672634467276, 0.05638672634467275, 3, 1, 1.351483666666667, 4.054451, 0.3806818235166093, 0.3806818235166093, 0.4391227214896667, 3.952104493407, 0.02146540180702206, 0.02146540180702206),
    @("MuSCs", "Il23a", "Il12rb1", "MuSCs", 2, 0.6666666666666666, 0.324919, 0.974757, 0.05638672634467276, 0.05638672634467275, 3, 1, 0.23571, 0.70713, 0.066394078474077, 0.066394078474077, 0.07658665749, 0.6892799174100001, 0.003743744733824508, 0.003743744733824507),
    @("MuSCs", "Il23a", "Il12rb1", "Resolving-Mac", 2, 0.6666666666666666, 0.324919, 0.974757, 0.05638672634467276, 0.05638672634467275, 3, 1, 0.7087373333333332, 2.126212, 0.1996349842044944, 0.1996349842044945, 0.2302822256093333, 2.072540030484, 0.0112567632231619, 0.0112567632231619),
    @("Resolving-Mac", "Il23a", "Il12rb1", "ECs", 3, 1, 1.065540666666667, 3.196622, 0.1849148556423401, 0.1849148556423401, 1, 0.3333333333333333, 0.2848286666666667, 0.854486, 0.08022967564521397, 0.08022967564521397, 0.3034965273657778, 2.731468746292, 0.01483565889016651, 0.01483565889016651),
    @("Resolving-Mac", "Il23a", "Il12rb1", "FAPs", 3, 1, 1.065540666666667, 3.196622, 0.1849148556423401, 0.1849148556423401, 3, 1, 0.9694063333333333, 2.908219, 0.2730594381596053, 0.2730594381596053, 1.032941870690889, 9.296476836218, 0.05049274658906192, 0.05049274658906191),
    @("Resolving-Mac", "Il23a", "Il12rb1", "Inflammatory-Mac", 3, 1, 1.065540666666667, 3.196622, 0.1849148556423401, 0.1849148556423401, 3, 1, 1.351483666666667, 4.054451, 0.3806818235166093, 0.3806818235166093, 1.440060807169111, 12.960547264522, 0.07039372444123661, 0.07039372444123661),
    @("Resolving-Mac", "Il23a", "Il12rb1", "MuSCs", 3, 1, 1.065540666666667, 3.196622, 0.1849148556423401, 0.1849148556423401, 3, 1, 0.23571, 0.70713, 0.066394078474077, 0.066394078474077, 0.25115859054, 2.26042731486, 0.01227725143654015, 0.01227725143654015),
    @("Resolving-Mac", "Il23a", "Il12rb1", "Resolving-Mac", 3, 1, 1.065540666666667, 3.196622, 0.1849148556423401, 0.1849148556423401, 3, 1, 0.7087373333333332, 2.126212, 0.1996349842044944, 0.1996349842044945, 0.7551884506515555, 6.796696055863999, 0.03691547428533494, 0.03691547428533494),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}
